# Updates Gugan's timesheet entries (rows 32-43) on sheet "28-04-2022"
# (the 18th / last tab) to reflect the revised task breakdown, and moves
# the sheet's scroll/selection position as recorded in the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(18)

# Row 32: Discussed Estimation / Project, 09:00 - 09:45
$ws.Range("B32").Value2 = "Discussed Estimation"
$ws.Range("C32").Value2 = "Project"
$ws.Range("D32").Value2 = 0.375
$ws.Range("E32").Value2 = 0.40625

# Row 33: Worked on HTML layout for User Creation Page / Project, 09:45 - 10:30
$ws.Range("B33").Value2 = "Worked on HTML layout for User Creation Page"
$ws.Range("C33").Value2 = "Project"
$ws.Range("D33").Value2 = 0.40625
$ws.Range("E33").Value2 = 0.4375

# Row 34: Morning Break / Lunch and Break, 10:30 - 11:00
$ws.Range("B34").Value2 = "Morning Break"
$ws.Range("C34").Value2 = "Lunch and Break"
$ws.Range("D34").Value2 = 0.4375
$ws.Range("E34").Value2 = 0.45833333333333331

# Row 35: Worked on HTML layout for User Creation Page / Project, 11:00 - 13:00
$ws.Range("B35").Value2 = "Worked on HTML layout for User Creation Page"
$ws.Range("C35").Value2 = "Project"
$ws.Range("D35").Value2 = 0.45833333333333331
$ws.Range("E35").Value2 = 0.54166666666666663

# Row 36: Customer Meeting / Project, 13:00 - 13:45
$ws.Range("B36").Value2 = "Customer Meeting"
$ws.Range("C36").Value2 = "Project"
$ws.Range("D36").Value2 = 0.54166666666666663
$ws.Range("E36").Value2 = 0.57291666666666663

# Row 37: Lunch Break / Lunch and Break, 13:45 - 14:15
$ws.Range("B37").Value2 = "Lunch Break"
$ws.Range("C37").Value2 = "Lunch and Break"
$ws.Range("D37").Value2 = 0.57291666666666663
$ws.Range("E37").Value2 = 0.59375

# Row 38: Team Meeting / Non Project, 14:30 - 15:15
$ws.Range("B38").Value2 = "Team Meeting"
$ws.Range("C38").Value2 = "Non Project"
$ws.Range("D38").Value2 = 0.60416666666666663
$ws.Range("E38").Value2 = 0.63541666666666663

# Row 39: Learned Angular <Topics> / Exploration, 15:30 - 17:00
$ws.Range("B39").Value2 = "Learned Angular <Topics>"
$ws.Range("C39").Value2 = "Exploration "
$ws.Range("D39").Value2 = 0.64583333333333337
$ws.Range("E39").Value2 = 0.70833333333333337

# Row 40: Evening Break / Lunch and Break, 17:10 - 17:30
$ws.Range("B40").Value2 = "Evening Break"
$ws.Range("C40").Value2 = "Lunch and Break"
$ws.Range("D40").Value2 = 0.71527777777777779
$ws.Range("E40").Value2 = 0.72916666666666663

# Row 41: Worked on HTML Layout for Wizard(skill)page / Project, 17:30 - 18:45
$ws.Range("B41").Value2 = "Worked on HTML Layout for Wizard(skill)page"
$ws.Range("C41").Value2 = "Project"
$ws.Range("D41").Value2 = 0.72916666666666663
$ws.Range("E41").Value2 = 0.78125

# Row 43: free-text note, no project/time entries
$ws.Range("B43").Value2 = "HTML layout wizard(In progress)"

# Restore the saved scroll position / selection for the sheet view
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C43").Select()
